$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2077.5
$ws.Range("I2").Value = 2160
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 2160
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -2047
$ws.Range("N2").Value = -1726

$ws.Range("H8").Value = 22044542
$ws.Range("I8").Value = 12400109
$ws.Range("J8").Value = 99200000
$ws.Range("K8").Value = 37200327
$ws.Range("L8").Value = 297600000
$ws.Range("M8").Value = -37200188
$ws.Range("N8").Value = -297600278

$ws.Range("H11").Value = 458
$ws.Range("I11").Value = 458
$ws.Range("K11").Value = 458
$ws.Range("M11").Value = -318

$ws.Range("H28").Value = 405.86667
$ws.Range("I28").Value = 420.66666
$ws.Range("K28").Value = 420.66666
$ws.Range("M28").Value = 64.33334000000002

$ws.Range("H38").Value = 461.27274
$ws.Range("I38").Value = 461.27274
$ws.Range("K38").Value = 1383.81822
$ws.Range("M38").Value = -1011.81822

$ws.Range("H39").Value = 2898
$ws.Range("I39").Value = 2664.4443
$ws.Range("K39").Value = 7993.3329
$ws.Range("M39").Value = -7697.3329

$ws.Range("H42").Value = 310.16666
$ws.Range("I42").Value = 233.57143
$ws.Range("K42").Value = 700.71429
$ws.Range("M42").Value = -470.71429

$ws.Range("H76").Value = 4393.0557
$ws.Range("I76").Value = 4393.0557
$ws.Range("K76").Value = 4393.0557
$ws.Range("M76").Value = -4078.0557

$ws.Range("H79").Value = 4393.0557
$ws.Range("I79").Value = 4393.0557
$ws.Range("K79").Value = 4393.0557
$ws.Range("M79").Value = -3301.0557

$ws.Range("H80").Value = 2215.353
$ws.Range("I80").Value = 2216.8333
$ws.Range("J80").Value = 2214.5454
$ws.Range("K80").Value = 6650.499899999999
$ws.Range("L80").Value = 6643.6362
$ws.Range("M80").Value = -5652.499899999999
$ws.Range("N80").Value = -8639.636200000001

$ws.Range("H83").Value = 2215.353
$ws.Range("I83").Value = 2216.8333
$ws.Range("J83").Value = 2214.5454
$ws.Range("K83").Value = 19951.4997
$ws.Range("L83").Value = 19930.9086
$ws.Range("M83").Value = -14959.4997
$ws.Range("N83").Value = -29914.9086

$ws.Range("H132").Value = 2083.6272
$ws.Range("I132").Value = 1006.4151
$ws.Range("K132").Value = 3019.2453
$ws.Range("M132").Value = -489.2453

$ws.Range("H137").Value = 3222
$ws.Range("I137").Value = 2576.8572
$ws.Range("J137").Value = 3974.6667
$ws.Range("K137").Value = 7730.571599999999
$ws.Range("L137").Value = 11924.0001
$ws.Range("M137").Value = -5180.571599999999
$ws.Range("N137").Value = -17024.0001

$ws.Range("H138").Value = 5293.197
$ws.Range("I138").Value = 4953.353
$ws.Range("J138").Value = 5424.5
$ws.Range("K138").Value = 14860.059
$ws.Range("L138").Value = 16273.5
$ws.Range("M138").Value = -9720.059000000001
$ws.Range("N138").Value = -26553.5

$ws.Range("H141").Value = 1222.3529
$ws.Range("I141").Value = 1249.0625
$ws.Range("J141").Value = 795
$ws.Range("K141").Value = 3747.1875
$ws.Range("L141").Value = 2385
$ws.Range("M141").Value = 1432.8125
$ws.Range("N141").Value = -12745

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1350.8334
$ws.Range("I21").Value = 1350.8334
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1350.8334
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -976.8334
$ws.Range("N21").ClearContents()

$ws.Range("H32").Value = 13444.796
$ws.Range("I32").Value = 11196.556
$ws.Range("K32").Value = 11196.556
$ws.Range("M32").Value = -10909.556

$ws.Range("H36").Value = 20000
$ws.Range("I36").Value = 20000
$ws.Range("K36").Value = 20000
$ws.Range("M36").Value = -19654

$ws.Range("H61").Value = 995
$ws.Range("I61").Value = 995
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 995
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -783
$ws.Range("N61").ClearContents()

$ws.Range("H102").Value = 922.4761999999999
$ws.Range("I102").Value = 935.4211
$ws.Range("K102").Value = 935.4211
$ws.Range("M102").Value = 686.5789

$ws.Range("H132").Value = 3131.4167
$ws.Range("I132").Value = 3109.7646
$ws.Range("K132").Value = 9329.293799999999
$ws.Range("M132").Value = -6799.293799999999

$ws.Range("H136").Value = 995
$ws.Range("I136").Value = 995
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2985
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -435
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18598.1
$ws.Range("I82").Value = 3667
$ws.Range("J82").Value = 24997.143
$ws.Range("K82").Value = 3667
$ws.Range("L82").Value = 24997.143
$ws.Range("M82").Value = -3284
$ws.Range("N82").Value = -25763.143

$ws.Range("H85").Value = 18598.1
$ws.Range("I85").Value = 3667
$ws.Range("J85").Value = 24997.143
$ws.Range("K85").Value = 3667
$ws.Range("L85").Value = 24997.143
$ws.Range("M85").Value = -2341
$ws.Range("N85").Value = -27649.143

$ws.Range("H86").Value = 5230.385
$ws.Range("I86").Value = 5842
$ws.Range("J86").Value = 4516.8335
$ws.Range("K86").Value = 5842
$ws.Range("L86").Value = 4516.8335
$ws.Range("M86").Value = -4719
$ws.Range("N86").Value = -6762.8335

$ws.Range("H89").Value = 5230.385
$ws.Range("I89").Value = 5842
$ws.Range("J89").Value = 4516.8335
$ws.Range("K89").Value = 29210
$ws.Range("L89").Value = 22584.1675
$ws.Range("M89").Value = -23594
$ws.Range("N89").Value = -33816.1675

$ws.Range("H134").Value = 76317
$ws.Range("I134").Value = 82299.164
$ws.Range("K134").Value = 246897.492
$ws.Range("M134").Value = -244362.492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32544.715
$ws.Range("I31").Value = 69467.664
$ws.Range("J31").Value = 4852.5
$ws.Range("K31").Value = 69467.664
$ws.Range("L31").Value = 4852.5
$ws.Range("M31").Value = -69172.664
$ws.Range("N31").Value = -5442.5

$ws.Range("H34").Value = 32544.715
$ws.Range("I34").Value = 69467.664
$ws.Range("J34").Value = 4852.5
$ws.Range("K34").Value = 69467.664
$ws.Range("L34").Value = 4852.5
$ws.Range("M34").Value = -69265.664
$ws.Range("N34").Value = -5256.5

$ws.Range("H97").Value = 29272.75
$ws.Range("J97").Value = 29272.75
$ws.Range("L97").Value = 29272.75
$ws.Range("N97").Value = -31254.75

$ws.Range("H132").Value = 3558.6099
$ws.Range("I132").Value = 3321.7026
$ws.Range("K132").Value = 9965.1078
$ws.Range("M132").Value = -7435.1078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9111.4
$ws.Range("I56").Value = 9111.4
$ws.Range("K56").Value = 9111.4
$ws.Range("M56").Value = -8581.4

$ws.Range("H92").Value = 1133
$ws.Range("I92").Value = 2999
$ws.Range("K92").Value = 8997
$ws.Range("M92").Value = -7749

$ws.Range("H107").Value = 677.2033699999999
$ws.Range("I107").Value = 415.23077
$ws.Range("K107").Value = 1245.69231
$ws.Range("M107").Value = 674.3076900000001

$ws.Range("H112").Value = 6062.25
$ws.Range("I112").Value = 6062.25
$ws.Range("K112").Value = 18186.75
$ws.Range("M112").Value = -17078.75

$ws.Range("H113").Value = 5424.857
$ws.Range("J113").Value = 5534.4614
$ws.Range("L113").Value = 16603.3842
$ws.Range("N113").Value = -20943.3842

$ws.Range("H131").Value = 2761.077
$ws.Range("J131").Value = 4157.5
$ws.Range("L131").Value = 12472.5
$ws.Range("N131").Value = -22552.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 44824.418
$ws.Range("J46").Value = 45978.6
$ws.Range("L46").Value = 45978.6
$ws.Range("N46").Value = -46290.6

$ws.Range("H132").Value = 8103.273
$ws.Range("I132").Value = 7690.4062
$ws.Range("J132").Value = 9204.25
$ws.Range("K132").Value = 23071.2186
$ws.Range("L132").Value = 27612.75
$ws.Range("M132").Value = -20541.2186
$ws.Range("N132").Value = -32672.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1261.8
$ws.Range("J22").Value = 1159.091
$ws.Range("L22").Value = 1159.091
$ws.Range("N22").Value = -1749.091

$ws.Range("H27").Value = 1261.8
$ws.Range("J27").Value = 1159.091
$ws.Range("L27").Value = 1159.091
$ws.Range("N27").Value = -1373.091

$ws.Range("H61").Value = 67331.34
$ws.Range("I61").Value = 51169.047
$ws.Range("K61").Value = 51169.047
$ws.Range("M61").Value = -50967.047

$ws.Range("H113").Value = 67331.34
$ws.Range("I113").Value = 51169.047
$ws.Range("K113").Value = 51169.047
$ws.Range("M113").Value = -48999.047

$ws.Range("H132").Value = 3963.3333
$ws.Range("I132").Value = 3963.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11889.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9359.999899999999
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 30271.047
$ws.Range("J136").Value = 9666.666999999999
$ws.Range("L136").Value = 29000.001
$ws.Range("N136").Value = -34100.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1817.1666
$ws.Range("I122").Value = 1680.6
$ws.Range("K122").Value = 5041.799999999999
$ws.Range("M122").Value = -2591.799999999999

$ws.Range("H132").Value = 7100.385
$ws.Range("I132").Value = 7100.385
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 21301.155
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -18771.155
$ws.Range("N132").ClearContents()
